$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Column D cells that hold numeric-looking text keep their exact
# string formatting (e.g. "1.00", "0.999") instead of being auto-converted
# to numbers by Excel, by forcing a Text number format before assignment.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "43.408.11"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").Value = "2.369.26"
$ws.Range("E3").Value = "  +2.42%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "310.02"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").Value = "103.82"
$ws.Range("E6").Value = "  +2.83%  "
$ws.Range("D7").Value = "0.510"
$ws.Range("E7").Value = "  -4.94%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "0.518"
$ws.Range("E9").Value = "  -1.32%  "
$ws.Range("D10").Value = "35.67"
$ws.Range("E10").Value = "  -1.04%  "
$ws.Range("D11").Value = "53.22"
$ws.Range("E11").Value = "  +1.91%  "
$ws.Range("E12").Value = "  -1.36%  "
$ws.Range("E14").Value = "  -3.96%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.744.15"
$ws.Range("E15").Value = "  +2.78%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "15.50"
$ws.Range("E16").Value = "  +3.42%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.372.37"
$ws.Range("E17").Value = "  +2.90%  "
$ws.Range("B18").Value = "Polygon"
$ws.Range("C18").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D18").Value = "0.808"
$ws.Range("E18").Value = "  -0.56%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "43.390.79"
$ws.Range("E19").Value = "  +0.78%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "6.31"
$ws.Range("E20").Value = "  +3.01%  "
$ws.Range("B21").Value = "InternetComputer(DFINITY)"
$ws.Range("C21").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D21").Value = "11.87"
$ws.Range("E21").Value = "  -5.74%  "
$ws.Range("B22").Value = "ShibaInu"
$ws.Range("C22").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D22").Value = "0.0X0913"
$ws.Range("D22").Characters(4,1).Text = [string][char]0x2083
$ws.Range("E22").Value = "  -0.81%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").Value = "68.10"
$ws.Range("E23").Value = "  -0.60%  "
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").Value = "239.91"
$ws.Range("E24").Value = "  -0.55%  "
$ws.Range("B25").Value = "ImmutableX"
$ws.Range("C25").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D25").Value = "2.04"
$ws.Range("E25").Value = "  +0.55%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").Value = "2.60"
$ws.Range("E26").Value = "  -1.18%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "25.77"
$ws.Range("E28").Value = "  +3.75%  "
$ws.Range("B29").Value = "LEO"
$ws.Range("C29").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D29").Value = "3.86"
$ws.Range("E29").Value = "  -2.60%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "2.31"
$ws.Range("E30").Value = "  +9.19%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").Value = "36.50"
$ws.Range("E31").Value = "  -2.33%  "
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").Value = "9.45"
$ws.Range("E32").Value = "  -2.45%  "
$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").Value = "161.54"
$ws.Range("E33").Value = "  -3.13%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "5.21"
$ws.Range("E34").Value = "  -2.51%  "
$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("B36").Value = "Celestia"
$ws.Range("C36").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D36").Value = "18.12"
$ws.Range("E36").Value = "  +0.67%  "
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").Value = "2.51"
$ws.Range("E37").Value = "  +4.67%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "4.66"
$ws.Range("E38").Value = "  +8.23%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").Value = "3.07"
$ws.Range("E39").Value = "  -3.01%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "0.0734"
$ws.Range("E40").Value = "  -1.24%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "1.91"
$ws.Range("E41").Value = "  +3.43%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "0.105"
$ws.Range("E42").Value = "  -1.91%  "
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").Value = "0.113"
$ws.Range("E43").Value = "  -2.49%  "
$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").Value = "2.61"
$ws.Range("E44").Value = "  +12.85%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.034.78"
$ws.Range("E45").Value = "  +2.83%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "19.60"
$ws.Range("E46").Value = "  +0.35%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "0.0289"
$ws.Range("E47").Value = "  -0.56%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "10.52"
$ws.Range("E48").Value = "  +6.89%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "3.09"
$ws.Range("E49").Value = "  +2.50%  "
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").Value = "57.64"
$ws.Range("E50").Value = "  +3.58%  "
$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D51").Value = "2.92"
$ws.Range("E51").Value = "  -1.81%  "
